$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the sample chess-pairing rows (rows 4-18, columns A-R),
# keeping cell styling but removing all values/content.
$ws.Range("A4:R18").ClearContents()

# Update the "push notification" username fields with their real values
$ws.Range("D1").Value = "White-UserName"
$ws.Range("M1").Value = "Black-UserName"
$ws.Range("D2").Value = "abcdg"
$ws.Range("M2").Value = "Test224"
$ws.Range("D3").Value = "SohamKale"
$ws.Range("M3").Value = "Harsh27"

# Move the active selection to K7
$ws.Range("K7").Select()
